# Update gh-pages output (generated at 456a3b4)
# Applies refreshed "want-to-go" counts (column F) across the four sheets
# and appends one newly-scraped local-life event row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (Exhibitions) - refresh column F counts
# ---------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F8").Value  = 37
$wsExpo.Range("F9").Value  = 717
$wsExpo.Range("F10").Value = 2660
$wsExpo.Range("F11").Value = 2660
$wsExpo.Range("F12").Value = 13
$wsExpo.Range("F13").Value = 1726
$wsExpo.Range("F14").Value = 599
$wsExpo.Range("F15").Value = 262
$wsExpo.Range("F16").Value = 678
$wsExpo.Range("F17").Value = 4940
$wsExpo.Range("F18").Value = 159
$wsExpo.Range("F22").Value = 852
$wsExpo.Range("F31").Value = 478
$wsExpo.Range("F33").Value = 798
$wsExpo.Range("F34").Value = 53
$wsExpo.Range("F37").Value = 1398
$wsExpo.Range("F38").Value = 1367

# ---------------------------------------------------------------
# Sheet "演出" (Performances) - refresh column F counts
# ---------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F11").Value = 134
$wsShow.Range("F18").Value = 256
$wsShow.Range("F22").Value = 1

# ---------------------------------------------------------------
# Sheet "本地生活" (Local life) - refresh column F counts and add
# the newly scraped event as a new row at the bottom of the table.
# ---------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 835
$wsLocal.Range("F4").Value = 238

# New row 8: copy the row-above's number-column formatting (bold,
# bordered, centered) onto A8 before filling in the values, and
# force column B's value to stay plain text so the literal date
# string isn't reinterpreted as a date serial number.
$wsLocal.Range("A7").Copy()
$wsLocal.Range("A8").PasteSpecial(-4122)

$wsLocal.Range("A8").Value = 7
$wsLocal.Range("B8").NumberFormat = "@"
$wsLocal.Range("B8").Value = "2024-10-28"
$wsLocal.Range("B8").Style = "Normal"
$wsLocal.Range("C8").Value = "北京·蜡笔小新：我们的恐龙日记x HAPPY ZOO 主题咖啡厅"
$wsLocal.Range("D8").Value = "王府井地铁站F1东口步行120米 北京王府井喜悦购物中心"
$wsLocal.Range("E8").Value = "2024.10.28 00:00-11.10 23:59"
$wsLocal.Range("F8").Value = 1
$wsLocal.Range("G8").Value = 10
$wsLocal.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=93723"
$wsLocal.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202410/zhKQZnHB1729477411932.png"

# ---------------------------------------------------------------
# Sheet "全部类型" (All types, aggregate) - refresh column F counts
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 835
$wsAll.Range("F7").Value  = 238
$wsAll.Range("F20").Value = 37
$wsAll.Range("F21").Value = 2660
$wsAll.Range("F23").Value = 1727
$wsAll.Range("F24").Value = 134
$wsAll.Range("F25").Value = 599
$wsAll.Range("F26").Value = 262
$wsAll.Range("F27").Value = 678
$wsAll.Range("F28").Value = 4940
$wsAll.Range("F32").Value = 852
$wsAll.Range("F41").Value = 478
$wsAll.Range("F44").Value = 256
$wsAll.Range("F46").Value = 798
$wsAll.Range("F47").Value = 53
$wsAll.Range("F50").Value = 1398
